# Updated cryptos list with GitHub Actions
# Applies the price/volume/row-order updates to the cryptos worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.135.86"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.650.33"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "19.75"
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.876.31"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "1.649.70"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "66.33"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").Value = "27.102.51"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "222.58"
$ws.Range("E19").Value = "  +3.58%  "
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").Value = "  +8.07%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.43"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "2.41"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "9.28"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "147.38"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "7.44"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "15.92"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "1.268.33"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").Value = "0.537"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "1.786.59"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "62.06"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "92.66"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -7.93%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.0975"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.63"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  -0.32%  "
